# BasecaseComparison.xlsx — fix an off-by-one in the data table: the
# "Emission" row needs a blank spacer row above it, and the "Tot"/"Oper"/
# "Capex" rows need a blank spacer row between them and "Emission". The
# numeric values are cleared out (to be re-entered later) and the blank
# spacer rows get a percentage number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above "Emission" (was row 2) - pushes Emission -> row 3,
# Tot -> row 4, Oper -> row 5, Capex -> row 6.
$ws.Rows.Item(2).EntireRow.Insert()

# Insert a second blank row below "Emission" (now row 3) - pushes
# Tot -> row 5, Oper -> row 6, Capex -> row 7.
$ws.Rows.Item(4).EntireRow.Insert()
# The second insert copies formatting down from row 3 ("Emission"); reset it
# back to the default/blank style so it matches the first spacer row.
$ws.Rows.Item(4).EntireRow.ClearFormats()

# Give both new spacer rows a percentage format.
$ws.Range("B2:D2").NumberFormat = "0.00%"
$ws.Range("B4:D4").NumberFormat = "0.00%"

# Clear the (now stale / to-be-recomputed) data values, keeping formatting.
$ws.Range("B3:D3").ClearContents()
$ws.Range("B5:D7").ClearContents()

# Match the author's final selection.
$ws.Range("D2:D7").Select()

# Re-point the chart series at their new rows.
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart

$s1 = $chart.SeriesCollection(1)
$s2 = $chart.SeriesCollection(2)
$s3 = $chart.SeriesCollection(3)
$s4 = $chart.SeriesCollection(4)

$s1.Formula = "=SERIES(Sheet1!`$A`$3,Sheet1!`$B`$1:`$D`$1,Sheet1!`$B`$3:`$D`$3,1)"
$s2.Formula = "=SERIES(Sheet1!`$A`$5,Sheet1!`$B`$1:`$D`$1,Sheet1!`$B`$5:`$D`$5,2)"
$s3.Formula = "=SERIES(Sheet1!`$A`$6,Sheet1!`$B`$1:`$D`$1,Sheet1!`$B`$6:`$D`$6,3)"
$s4.Formula = "=SERIES(Sheet1!`$A`$7,Sheet1!`$B`$1:`$D`$1,Sheet1!`$B`$7:`$D`$7,4)"

# The chart anchor needs to move down (two new rows above it) and grow to
# match the author's resize.
$chartObj.Left = 355.625
$chartObj.Top = 83
$chartObj.Width = 560.9375
$chartObj.Height = 363
